$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Led" quantity 4 -> 2
$ws.Range("B7").Value = 2

# "Motor Servo" -> "Buzzer" (quantity unchanged at 1)
$ws.Range("A9").Value = "Buzzer"

# Normalize the leftover duplicate style on A4:B4 (copy the plain row style
# used elsewhere in the table) so the redundant style entry is dropped.
$ws.Range("A2:B2").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)

# Match the final selection state left by the author
$ws.Range("B9").Select()
